$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(3046, "[0, 5, -7]", 15.28572233709476),
  @(3047, "[-1, 6, -7]", 14.97555531673779),
  @(3048, "[7, -2, -6]", 14.06195343883236),
  @(3049, "[6, -1, -4]", 14.43407948864752),
  @(3050, "[4, 2, -6]", 14.26554181851618),
  @(3051, "[0, 5, -5]", 15.28884288230845),
  @(3052, "[1, 4, -6]", 15.20394382332994),
  @(3053, "[-4, 6, -4]", 13.93494283265519),
  @(3054, "[0, 5, -4]", 15.30325754584289),
  @(3055, "[-1, 5, -4]", 15.04361706071651),
  @(3056, "[-1, 6, -4]", 15.04131883169828),
  @(3057, "[0, 6, -2]", 15.85964945753043),
  @(3058, "[5, 2, -6]", 14.10379931732365),
  @(3059, "[7, -2, -4]", 14.16188233000509),
  @(3060, "[-1, 5, -3]", 15.00768879131732),
  @(3061, "[7, -1, -5]", 14.33638116794155),
  @(3062, "[0, 6, -6]", 15.3058721219751),
  @(3063, "[0, 6, -7]", 15.29083894892729),
  @(3064, "[0, 7, -7]", 15.30567959749544),
  @(3065, "[7, -1, -4]", 14.40713168736626),
  @(3066, "[0, 6, -1]", 18.20096342438711),
  @(3067, "[7, -1, -3]", 14.574244612021),
  @(3068, "[4, 3, -4]", 14.05130368380875),
  @(3069, "[2, 4, -7]", 14.80809407701546),
  @(3070, "[-1, 6, -3]", 15.00815456990874),
  @(3071, "[1, 5, -5]", 15.1646689388156),
  @(3072, "[-3, 6, -4]", 14.21226413693196),
  @(3073, "[7, -2, -7]", 14.05717189314715),
  @(3074, "[7, 0, -4]", 14.7503870503488),
  @(3075, "[6, -2, -7]", 14.06780589097803),
  @(3076, "[-3, 7, -5]", 14.16629096421105),
  @(3077, "[5, 1, -6]", 14.08681626842426),
  @(3078, "[6, 1, -6]", 14.37187711524859),
  @(3079, "[6, 1, -7]", 14.37086333476918),
  @(3080, "[6, 0, -7]", 14.66662159871884),
  @(3081, "[6, -1, -7]", 14.35013278483938),
  @(3082, "[7, -1, -7]", 14.32892505968756),
  @(3083, "[1, 5, -7]", 15.15443983808793),
  @(3084, "[3, 3, -5]", 14.22336862000474),
  @(3085, "[1, 5, -6]", 15.16977113423794),
  @(3086, "[-1, 7, -6]", 15.00169111896097),
  @(3087, "[0, 6, -5]", 15.29647891018559)
)

$startRow = 3048
$endRow = $startRow + $data.Count - 1

# Replicate the formatting of column A (bold font, thin border, centered)
# from the last existing data row onto the whole block of new rows in one shot.
$ws.Cells.Item(3047, 1).Copy()
$ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 1)).PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $startRow + $i
  $item = $data[$i]
  $ws.Cells.Item($row, 1).Value = $item[0]
  $ws.Cells.Item($row, 2).Value = $item[1]
  $ws.Cells.Item($row, 3).Value = $item[2]
}

Write-Host "Added " $data.Count " rows"
